# Fixed non-adaptive test, new EDT2 item bank
# Adds four new rows (33-36) to the EDT sheet: ANGER/HAPPY/SAD/FEAR item-bank
# keys with their DE/EN/RU/NL translations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (key) -------------------------------------------------
$ws.Range("A33").Value = "ANGER"
$ws.Range("A34").Value = "HAPPY"
$ws.Range("A35").Value = "SAD"
$ws.Range("A36").Value = "FEAR"

# --- Column B (DE) ----------------------------------------------------
$ws.Range("B33").Value = "am wütendsten"
$ws.Range("B34").Value = "am fröhlichsten"
$ws.Range("B35").Value = "am traurigsten"
$ws.Range("B36").Value = "am ängstlichsten"

# --- Column C (EN) -- entered in this specific order to match the
# original authoring order of the shared-string table -----------------
$ws.Range("C34").Value = "happiest"
$ws.Range("C33").Value = "angriest"
$ws.Range("C35").Value = "saddest"
$ws.Range("C36").Value = "most fearful"

# --- Column D (RU) -- reuses already-existing shared strings ----------
$ws.Range("D33").Value = "агрессивнее"
$ws.Range("D34").Value = "веселее"
$ws.Range("D35").Value = "печальнее"
$ws.Range("D36").Value = "более пугающей"

# --- Column E (NL) -- reuses already-existing shared strings ----------
$ws.Range("E33").Value = "boos"
$ws.Range("E34").Value = "vrolijk"
$ws.Range("E35").Value = "verdrietig"
$ws.Range("E36").Value = "angstig"

# Row 33 is taller, matching the other "header-ish" rows in the sheet
$ws.Rows(33).RowHeight = 16

# E33 picks up the wrap-text style used elsewhere in column E (e.g. E7)
$ws.Range("E7").Copy()
$ws.Range("E33").PasteSpecial(-4122)

# D36 picks up the "Calibri (Textkörper)" font style used at D14
$ws.Range("D14").Copy()
$ws.Range("D36").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update the selection / active cell as recorded in the saved view
$ws.Range("C38").Select() | Out-Null
